$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for the new columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy style from an existing header cell (AC1) so the new headers match formatting
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in team record data for each data row (rows 2 through 42)
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 92   # AD = column 30
    $ws.Cells.Item($r, 31).Value = 70   # AE = column 31
    $ws.Cells.Item($r, 32).Value = 0    # AF = column 32
}
